# Jenny's FAQ comments: add a new bullet " Why is the asymmetry not there?"
# right before " Assumptions on Ehat? ", and carry the "_GoBack" bookmark
# (which marks the most-recent edit location) along with the newly typed
# bullet instead of leaving it on the old last bullet.

$d = $word.ActiveDocument

# --- Step 1: drop the existing _GoBack bookmark (it currently sits at the
# end of the last bullet, " Do you really believe ... interest rate?").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: locate the " Assumptions on Ehat? " bullet and insert a brand
# new list paragraph right before it containing the new question. Because
# the insertion point sits at the start of an existing ListParagraph/numPr
# paragraph, the freshly split-off paragraph inherits the same list
# formatting automatically.
$target = $d.Range(0, $d.Content.End)
$found = $target.Find.Execute(" Assumptions on Ehat? ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetParaIndex = $target.Paragraphs(1).Index
$insertionPoint = $d.Range($target.Start, $target.Start)
$insertionPoint.InsertBefore(" Why is the asymmetry not there?`r")

# --- Step 3: re-home the _GoBack bookmark onto the end of the bullet we
# just typed, exactly like Word does when you finish typing new text.
# A zero-length Range placed directly at a paragraph-end position cannot
# be handed straight to Bookmarks.Add, so we bookmark the final character
# of the new bullet, delete that single character (collapsing the
# bookmark to zero width at that spot), and retype the character in front
# of the now-collapsed bookmark so the visible text is unchanged.
$newPara = $d.Paragraphs($targetParaIndex)
$bulletEnd = $newPara.Range.End - 1
$lastCharRange = $d.Range($bulletEnd - 1, $bulletEnd)
$savedChar = $lastCharRange.Text

$d.Bookmarks.Add("_GoBack", $lastCharRange)
$tempBookmark = $d.Bookmarks("_GoBack")
$tempBookmark.Range.Delete()

$collapsedBookmark = $d.Bookmarks("_GoBack")
$retypePoint = $d.Range($collapsedBookmark.Start, $collapsedBookmark.Start)
$retypePoint.InsertBefore($savedChar)

Write-Output "Inserted new bullet and relocated _GoBack bookmark."
